# Selumetinib_SRB_data.xlsx - "updating daily entry + obj summaries"
#
# For both S4 and S5 sheets:
#  - highlight the per-day replicate block (B.. : K..) with the existing
#    yellow "summary block" style (same style already used on S1/S2/S3)
#  - add a row-16 AVERAGE() summary row (shared formula across C:K)
#  - add an 18-27 block: col B = the row-16 averages transposed down the
#    column (plain values), col C = row-18 is the literal baseline 100,
#    rows 19-27 are a shared "% of day-1 average" formula
#
# S4 becomes the active/visible tab (was S5).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("S1")
$ws2 = $wb.Worksheets.Item("S2")
$ws4 = $wb.Worksheets.Item("S4")
$ws5 = $wb.Worksheets.Item("S5")

# ---------------------------------------------------------------------
# S4 ("S4" / sheet4.xml)
# ---------------------------------------------------------------------

# Highlight B8:K12 with the yellow block style already used elsewhere in
# this workbook (S2/S3 already carry it) - copy format only, so it reuses
# the existing style index instead of minting a new one.
$ws2.Range("B8").Copy() | Out-Null
$ws4.Range("B8:K12").PasteSpecial(-4122) | Out-Null

# Row 16: per-column averages of the highlighted block (B8:B12 .. K8:K12)
$ws4.Range("B16").Formula = "=AVERAGE(B8:B12)"
$ws4.Range("C16:K16").Formula = "=AVERAGE(C8:C12)"

# Rows 18-27: col B = transposed row-16 averages (plain values); col C =
# percent-of-day-1 with the existing percent-block style.
$ws1.Range("C25").Copy() | Out-Null
$ws4.Range("C18:C27").PasteSpecial(-4122) | Out-Null

$ws4.Range("B18").Value = $ws4.Range("B16").Value()
$ws4.Range("B19").Value = $ws4.Range("C16").Value()
$ws4.Range("B20").Value = $ws4.Range("D16").Value()
$ws4.Range("B21").Value = $ws4.Range("E16").Value()
$ws4.Range("B22").Value = $ws4.Range("F16").Value()
$ws4.Range("B23").Value = $ws4.Range("G16").Value()
$ws4.Range("B24").Value = $ws4.Range("H16").Value()
$ws4.Range("B25").Value = $ws4.Range("I16").Value()
$ws4.Range("B26").Value = $ws4.Range("J16").Value()
$ws4.Range("B27").Value = $ws4.Range("K16").Value()

$ws4.Range("C18").Value = 100
$ws4.Range("C19").Formula = "=(B19/0.1738*100)"
$ws4.Range("C20:C27").Formula = "=(B20/0.1738*100)"

# View state: S4 becomes the active sheet/tab.
$ws4.Activate()
$ws4.Range("F36").Select() | Out-Null

# ---------------------------------------------------------------------
# S5 ("S5" / sheet5.xml)
# ---------------------------------------------------------------------

# Highlight B11:K13 with the same yellow block style.
$ws2.Range("B8").Copy() | Out-Null
$ws5.Range("B11:K13").PasteSpecial(-4122) | Out-Null

# Row 16: per-column averages of the highlighted block (B11:B13 .. K11:K13)
$ws5.Range("B16").Formula = "=AVERAGE(B11:B13)"
$ws5.Range("C16:K16").Formula = "=AVERAGE(C11:C13)"

# Rows 18-27: col B = transposed row-16 averages (plain values); col C =
# percent-of-day-1 with the existing percent-block style.
$ws1.Range("C25").Copy() | Out-Null
$ws5.Range("C18:C27").PasteSpecial(-4122) | Out-Null

$ws5.Range("B18").Value = $ws5.Range("B16").Value()
$ws5.Range("B19").Value = $ws5.Range("C16").Value()
$ws5.Range("B20").Value = $ws5.Range("D16").Value()
$ws5.Range("B21").Value = $ws5.Range("E16").Value()
$ws5.Range("B22").Value = $ws5.Range("F16").Value()
$ws5.Range("B23").Value = $ws5.Range("G16").Value()
$ws5.Range("B24").Value = $ws5.Range("H16").Value()
$ws5.Range("B25").Value = $ws5.Range("I16").Value()
$ws5.Range("B26").Value = $ws5.Range("J16").Value()
$ws5.Range("B27").Value = $ws5.Range("K16").Value()

$ws5.Range("C18").Value = 100
$ws5.Range("C19").Formula = "=B19/0.04333*100"
$ws5.Range("C20:C27").Formula = "=B20/0.04333*100"

# View state: S5 keeps a plain (non-active) view, selection moves to K25:K26.
$ws5.Range("K25:K26").Select() | Out-Null

# Re-activate S4 last so it ends up the workbook's active/visible tab.
$ws4.Activate()
